$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "160"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "410599.00"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "314"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1054917.92"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "151"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "415137.26"

$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "8"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "32700.00"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "233"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "594019.00"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "471"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1776312.70"

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "9"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21000.00"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "292"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "789606.74"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "535"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1942644.11"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "352"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1140381.16"

$ws.Range("C72").NumberFormat = "@"
$ws.Range("C72").Value = "348"
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = "858635.70"

$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "848"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "2791581.34"

$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "485"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "1510920.87"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "200"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "465071.00"

$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "4"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "9500.00"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "472"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "1574000.50"

$ws.Range("C87").NumberFormat = "@"
$ws.Range("C87").Value = "174"
$ws.Range("D87").NumberFormat = "@"
$ws.Range("D87").Value = "487976.09"

$ws.Range("C89").NumberFormat = "@"
$ws.Range("C89").Value = "7"
$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "22670.00"

